$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing account 004216504 / WANDER (Excel row 7),
# shifting subsequent rows up.
$ws.Rows.Item(7).Delete()
